$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44363
$ws.Range("Q3").Value = '$/caja 15 kilos empedrada'
$ws.Range("S3").Value = 633
$ws.Range("T3").Value = 15

# Row 4
$ws.Range("D4").Value = 44316
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9500
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 528

# Row 5
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 11000
$ws.Range("P5").Value = 10500
$ws.Range("S5").Value = 583

# Row 6
$ws.Range("D6").Value = 44299
$ws.Range("L6").Value = 'Segunda'
$ws.Range("N6").Value = 9000
$ws.Range("O6").Value = 9000
$ws.Range("P6").Value = 9000
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región del Maule'
$ws.Range("S6").Value = 500

# Row 7
$ws.Range("D7").Value = 44272
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 9500
$ws.Range("Q7").Value = '$/caja 15 kilos granel'
$ws.Range("S7").Value = 633
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("S8").Value = 533

# Row 9
$ws.Range("D9").Value = 44307
$ws.Range("L9").Value = 'Primera'
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("Q9").Value = '$/bandeja 18 kilos granel'
$ws.Range("S9").Value = 556
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44307
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 8000
$ws.Range("O10").Value = 8000
$ws.Range("P10").Value = 8000
$ws.Range("Q10").Value = '$/bandeja 18 kilos granel'
$ws.Range("S10").Value = 444
$ws.Range("T10").Value = 18
